# Apply the vocabulary.xlsx update:
# 1. Update the "dct:modified" timestamp in B21.
# 2. Re-order several mapping-definition cells in row 23 (C,D,G and V,W,X,Y,Z).
# 3. Clear the stray "var" value in D25.
# 4. Move the "vocab:1002" broader-reference from column G to column C for
#    rows 27-29, replacing the old ad-hoc labels.
# 5. Move the "owl:ObjectProperty" rdf:type value from V30 to D30.
# 6. Remove the now-unused trailing column AO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: refreshed modification timestamp ---
$ws.Range("B21").Value = "2023-09-13T15:10:53+00:00"

# --- Row 23: header / mapping-definition row ---
$ws.Range("C23").Value = 'skos:broader(separator=",")'
$ws.Range("D23").Value = "rdf:type"
$ws.Range("G23").Value = 'skos:altLabel(separator=",")'
$ws.Range("V23").Value = "dct:modified^^xsd:date"
$ws.Range("W23").Value = "dct:created^^xsd:date"
$ws.Range("X23").Value = 'dct:creator(separator=",")'
$ws.Range("Y23").Value = 'dct:contributor(separator=",")'
$ws.Range("Z23").Value = ""

# --- Row 25: clear the leftover "var" notation value ---
$ws.Range("D25").Value = ""

# --- Row 27 (emerging): replace "new" label with broader reference ---
$ws.Range("C27").Value = "vocab:1002"
$ws.Range("G27").Value = ""

# --- Row 28 (developing): replace "intermediate" label with broader reference ---
$ws.Range("C28").Value = "vocab:1002"
$ws.Range("G28").Value = ""

# --- Row 29 (mature): add broader reference (was only in G) ---
$ws.Range("C29").Value = "vocab:1002"
$ws.Range("G29").Value = ""

# --- Row 30 (hasMaturityLevel): move rdf:type value from V to D ---
$ws.Range("D30").Value = "owl:ObjectProperty"
$ws.Range("V30").Value = ""

# --- Remove the now-empty trailing column AO (A1:AO30 -> A1:AN30) ---
$ws.Range("AO1:AO30").EntireColumn.Delete()
